$d = $word.ActiveDocument

$d.Content.Find.Execute("63×33=", $true, $false, $false, $false, $false, $true, 1, $false, "83×27=", 2) | Out-Null
$d.Content.Find.Execute("37×72=", $true, $false, $false, $false, $false, $true, 1, $false, "47×19=", 2) | Out-Null
$d.Content.Find.Execute("39×20=", $true, $false, $false, $false, $false, $true, 1, $false, "80×19=", 2) | Out-Null
$d.Content.Find.Execute("73×93=", $true, $false, $false, $false, $false, $true, 1, $false, "31×46=", 2) | Out-Null
$d.Content.Find.Execute("47×72=", $true, $false, $false, $false, $false, $true, 1, $false, "21×83=", 2) | Out-Null
$d.Content.Find.Execute("81×83=", $true, $false, $false, $false, $false, $true, 1, $false, "49×55=", 2) | Out-Null
$d.Content.Find.Execute("98×39=", $true, $false, $false, $false, $false, $true, 1, $false, "96×45=", 2) | Out-Null
$d.Content.Find.Execute("87×38=", $true, $false, $false, $false, $false, $true, 1, $false, "27×27=", 2) | Out-Null
$d.Content.Find.Execute("37×18=", $true, $false, $false, $false, $false, $true, 1, $false, "55×97=", 2) | Out-Null
$d.Content.Find.Execute("14×14=", $true, $false, $false, $false, $false, $true, 1, $false, "47×88=", 2) | Out-Null
$d.Content.Find.Execute("38×29=", $true, $false, $false, $false, $false, $true, 1, $false, "51×13=", 2) | Out-Null
$d.Content.Find.Execute("35×42=", $true, $false, $false, $false, $false, $true, 1, $false, "40×77=", 2) | Out-Null
$d.Content.Find.Execute("68×56=", $true, $false, $false, $false, $false, $true, 1, $false, "25×14=", 2) | Out-Null
$d.Content.Find.Execute("86×49=", $true, $false, $false, $false, $false, $true, 1, $false, "17×86=", 2) | Out-Null
$d.Content.Find.Execute("24×61=", $true, $false, $false, $false, $false, $true, 1, $false, "54×43=", 2) | Out-Null
$d.Content.Find.Execute("42×68=", $true, $false, $false, $false, $false, $true, 1, $false, "77×43=", 2) | Out-Null
$d.Content.Find.Execute("19×70=", $true, $false, $false, $false, $false, $true, 1, $false, "63×50=", 2) | Out-Null
$d.Content.Find.Execute("93×76=", $true, $false, $false, $false, $false, $true, 1, $false, "46×62=", 2) | Out-Null
$d.Content.Find.Execute("12×56=", $true, $false, $false, $false, $false, $true, 1, $false, "15×78=", 2) | Out-Null
$d.Content.Find.Execute("39×77=", $true, $false, $false, $false, $false, $true, 1, $false, "63×72=", 2) | Out-Null
$d.Content.Find.Execute("33×19=", $true, $false, $false, $false, $false, $true, 1, $false, "97×60=", 2) | Out-Null
$d.Content.Find.Execute("29×26=", $true, $false, $false, $false, $false, $true, 1, $false, "79×49=", 2) | Out-Null
$d.Content.Find.Execute("39×96=", $true, $false, $false, $false, $false, $true, 1, $false, "95×14=", 2) | Out-Null
$d.Content.Find.Execute("55×70=", $true, $false, $false, $false, $false, $true, 1, $false, "35×31=", 2) | Out-Null
$d.Content.Find.Execute("76×72=", $true, $false, $false, $false, $false, $true, 1, $false, "49×66=", 2) | Out-Null

Write-Output "Replacements complete"
